$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (Metadata sheet, first sheet) ---
$ws1 = $wb.Worksheets.Item("Metadata")

# URL value (B2): fr/medication -> ig/fhir/medication
$ws1.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-medication-history-source-type"

# Name value (B4): FrMedicationHistorySourceType -> FRMedicationHistorySourceType
$ws1.Range("B4").Value = "FRMedicationHistorySourceType"

# Date value (B8): updated timestamp
$ws1.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction value (B11): was empty -> now FRANCE
$ws1.Range("B11").Value = "FRANCE"

# --- Sheet "Include #0" (second sheet) ---
$ws2 = $wb.Worksheets.Item("Include #0")

# System URI value (B4): fr/medication -> ig/fhir/medication
$ws2.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-history-source-type"
